$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: report generated timestamp ---
$ws.Range("D5").Value = "Report Generated On: 08/18/2025 09:47 PM"

# --- Summary box: Total Billed Amount / Total Line Items ---
$ws.Range("C8").Value = 0
$ws.Range("C9").Value = 4

# --- Insert a new detail row before the TOTAL row (old row 19) ---
# This pushes the old row 19 (TOTAL) down to row 20, and shifts the
# "A19:G19" merged cell to "A20:G20" automatically.
$ws.Rows.Item(19).Insert()

# Copy the banding/format (styles 12/13/14) from row 17 onto the newly
# inserted row 19 so it continues the alternating row-stripe pattern.
$ws.Range("A17:H17").Copy()
$ws.Range("A19").PasteSpecial(-4122)

# --- Row 16: PLA-DLOC/Inst -> POL-40-2/Inst, counts zeroed out ---
$ws.Range("B16").Value = "POL-40-2"
$ws.Range("D16").Value = "Pole,40ft,Class 2"
$ws.Range("F16").Value = 0
$ws.Range("H16").Value = 0

# --- Row 17: PLA-DLOC/Rem -> PLA-DLOC/Inst, billed amount zeroed ---
$ws.Range("C17").Value = "Inst"
$ws.Range("H17").Value = 0

# --- Row 18: PLA-BACK/Inst -> PLA-DLOC/Rem, qty + billed amount updated ---
$ws.Range("B18").Value = "PLA-DLOC"
$ws.Range("C18").Value = "Rem"
$ws.Range("D18").Value = "PLA,Difficult Location"
$ws.Range("F18").Value = 4
$ws.Range("H18").Value = 0

# --- Row 19 (newly inserted): PLA-BACK / Inst line item ---
$ws.Range("A19").Value = "Point 11"
$ws.Range("B19").Value = "PLA-BACK"
$ws.Range("C19").Value = "Inst"
$ws.Range("D19").Value = "Difficult Location Equip Adder-Backyard"
$ws.Range("E19").Value = "EA"
$ws.Range("F19").Value = 18
$ws.Range("G19").Value = ""
$ws.Range("H19").Value = 0

# --- Row 20 (was row 19): TOTAL, now zero ---
$ws.Range("H20").Value = 0
